$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates derived from the crypto price-ticker refresh.
# Numeric-looking strings (e.g. "0.9998") must be written to a
# Text-formatted cell so Excel does not silently coerce them to
# doubles; the NumberFormat is then reset back to the default
# "Normal" style so no stray formatting is left behind.

$ws.Range('D2').Value = '29.339.01'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').Value = '1.842.23'
$ws.Range('E3').Value = '  -0.28%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9983'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6292'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.56%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2945'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.42%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07686'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.44%  '

$ws.Range('D12').Value = '1.863.67'
$ws.Range('E12').Value = '  -6.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.976'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6779'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.39%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001026'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.14%  '

$ws.Range('D17').Value = '2.112.55'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.145'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.10%  '

$ws.Range('D19').Value = '29.370.11'
$ws.Range('E19').Value = '  -0.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.90%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.94%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.457'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.29%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.31%  '

$ws.Range('E26').Value = '  +0.42%  '

$ws.Range('E27').Value = '  -1.33%  '

$ws.Range('E28').Value = '  -0.33%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.459'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.269'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.05%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05648'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.83%  '

$ws.Range('E32').Value = '  -0.25%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.018'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.04%  '

$ws.Range('E34').Value = '  -1.95%  '

$ws.Range('E35').Value = '  -0.47%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7103'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.95%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.589'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.42%  '

$ws.Range('D38').Value = '1.240.51'
$ws.Range('E38').Value = '  -0.77%  '

$ws.Range('E39').Value = '  +0.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.770'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.81%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.235'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.33%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9010'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9992'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.05%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.67%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000118'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.05%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.087'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.28%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.3993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.77%  '

$ws.Range('E49').Value = '  -1.38%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.935'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.33%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1120'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.10%  '
